$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new numeric-looking value must stay as text
$textCells = @("D4", "D5", "D6", "D10", "D13", "D14", "D18", "D21", "D22", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D35", "D37", "D39", "D40", "D41", "D42", "D45", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row value updates
$ws.Range("D2").Value = "68.820.52"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.851.44"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "601.74"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "172.19"
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("D7").Value = "3.851.94"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("E11").Value = "  +3.74%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "0.0000287"
$ws.Range("E13").Value = "  +15.72%  "
$ws.Range("D14").Value = "37.28"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "4.497.65"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "3.899.98"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").Value = "68.819.90"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "18.37"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "11.26"
$ws.Range("E21").Value = "  +4.83%  "
$ws.Range("D22").Value = "474.69"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").Value = "83.82"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").Value = "10.41"
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "2.95"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "4.001.96"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "7.77"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "31.41"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "9.36"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "3.817.56"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "4.00"
$ws.Range("E37").Value = "  +22.10%  "
$ws.Range("E38").Value = "  -0.44%  "

# Rows 39/40 swap places: Mantle <-> Kaspa
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "1.02"
$ws.Range("E40").Value = "  +1.08%  "

$ws.Range("D41").Value = "5.95"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "0.000301"
$ws.Range("E45").Value = "  +10.43%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "420.70"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").Value = "46.50"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").Value = "142.17"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("E51").Value = "  +0.19%  "
